$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "preguntas:"
$ws.Range("A9").Value = "1.- Todo estos campos se exportan o algunos los gestiono yo. Ej: tipo, nivel, accion, estado, fecha citacion, incidencia atendida, fecha atendida,  hora atendida, fecha cerrada, hora cerrada, usuario cerrada, motivo de erro"
$ws.Range("A10").Value = "2.- Las cabeceras son nombres definitivos, necesito que lo sean respentando todos los caracteres, espacios, tildes, …"
$ws.Range("A11").Value = "3.- Entiendo que el codigo incidencia será unico. Lo necesitaría para no introducir en la base de datos incidencias repetidas"

$ws.Range("A1").Select()
